$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5679
$ws.Range("E2").Value = 333
$ws.Range("F2").Value = 333
$ws.Range("G2").Value = 217
$ws.Range("H2").Value = 179
$ws.Range("I2").Value = 153
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 5218
$ws.Range("L2").Value = 3841
$ws.Range("M2").Value = 1377
$ws.Range("N2").Value = 1171
$ws.Range("O2").Value = 206
$ws.Range("P2").Value = 205
$ws.Range("Q2").Value = 244
$ws.Range("R2").Value = -156
$ws.Range("S2").Value = -23
$ws.Range("T2").Value = 104
$ws.Range("U2").Value = 140
$ws.Range("V2").Value = 2987
$ws.Range("W2").Value = 5.86
$ws.Range("X2").Value = 3.15
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").Value = 278.9
$ws.Range("AB2").Value = 470.75
$ws.Range("AC2").Value = 443
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").Value = 3304
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 35432000

# Row 3
$ws.Range("D3").Value = 4533
$ws.Range("E3").Value = 354
$ws.Range("F3").Value = 430
$ws.Range("G3").Value = 243
$ws.Range("H3").Value = 238
$ws.Range("I3").Value = 198
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 5203
$ws.Range("L3").Value = 3500
$ws.Range("M3").Value = 1702
$ws.Range("N3").Value = 1459
$ws.Range("O3").Value = 243
$ws.Range("P3").Value = 219
$ws.Range("Q3").Value = 472
$ws.Range("R3").Value = -243
$ws.Range("S3").Value = -164
$ws.Range("T3").Value = 231
$ws.Range("U3").Value = 241
$ws.Range("V3").Value = 2597
$ws.Range("W3").Value = 7.81
$ws.Range("X3").Value = 5.25
$ws.Range("Y3").Value = 15.06
$ws.Range("Z3").Value = 4.57
$ws.Range("AA3").Value = 205.6
$ws.Range("AB3").Value = 564.66
$ws.Range("AC3").Value = 554
$ws.Range("AD3").Value = 5.55
$ws.Range("AE3").Value = 3810
$ws.Range("AF3").Value = 0.8100000000000001
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 3.25
$ws.Range("AI3").Value = 19.33
$ws.Range("AJ3").Value = 38292720

# Row 4
$ws.Range("D4").Value = 4274
$ws.Range("E4").Value = 197
$ws.Range("F4").Value = 197
$ws.Range("G4").Value = 143
$ws.Range("H4").Value = 197
$ws.Range("I4").Value = 141
$ws.Range("J4").Value = 55
$ws.Range("K4").Value = 3658
$ws.Range("L4").Value = 2075
$ws.Range("M4").Value = 1584
$ws.Range("N4").Value = 1584
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 220
$ws.Range("Q4").Value = 163
$ws.Range("R4").Value = 109
$ws.Range("S4").Value = -270
$ws.Range("T4").Value = 90
$ws.Range("U4").Value = 73
$ws.Range("V4").Value = 1751
$ws.Range("W4").Value = 4.62
$ws.Range("X4").Value = 4.61
$ws.Range("Y4").Value = 9.300000000000001
$ws.Range("Z4").Value = 4.44
$ws.Range("AA4").Value = 130.98
$ws.Range("AB4").Value = 620.67
$ws.Range("AC4").Value = 369
$ws.Range("AD4").Value = 8.460000000000001
$ws.Range("AE4").Value = 4133
$ws.Range("AF4").Value = 0.76
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 3.2
$ws.Range("AI4").Value = 27.09
$ws.Range("AJ4").Value = 38322720

# Row 5
$ws.Range("D5").Value = 2518
$ws.Range("E5").Value = 71
$ws.Range("F5").Value = 71
$ws.Range("G5").Value = 24
$ws.Range("H5").Value = 19
$ws.Range("I5").Value = 19
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 3490
$ws.Range("L5").Value = 1920
$ws.Range("M5").Value = 1570
$ws.Range("N5").Value = 1570
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 220
$ws.Range("Q5").Value = 241
$ws.Range("R5").Value = -123
$ws.Range("S5").Value = -123
$ws.Range("T5").Value = 14
$ws.Range("U5").Value = 228
$ws.Range("V5").Value = 1675
$ws.Range("W5").Value = 2.83
$ws.Range("X5").Value = 0.74
$ws.Range("Y5").Value = 1.18
$ws.Range("Z5").Value = 0.52
$ws.Range("AA5").Value = 122.27
$ws.Range("AB5").Value = 614.29
$ws.Range("AC5").Value = 49
$ws.Range("AD5").Value = 57.93
$ws.Range("AE5").Value = 4092
$ws.Range("AF5").Value = 0.6899999999999999
$ws.Range("AG5").Value = 30
$ws.Range("AH5").Value = 1.06
$ws.Range("AI5").Value = 61.64
$ws.Range("AJ5").Value = 38372720

# Row 6
$ws.Range("D6").Value = 2182
$ws.Range("E6").Value = 84
$ws.Range("F6").Value = 84
$ws.Range("G6").Value = 21
$ws.Range("H6").Value = 18
$ws.Range("I6").Value = 18
$ws.Range("K6").Value = 3683
$ws.Range("L6").Value = 2099
$ws.Range("M6").Value = 1584
$ws.Range("N6").Value = 1554
$ws.Range("P6").Value = 220
$ws.Range("Q6").Value = 210
$ws.Range("R6").Value = -185
$ws.Range("S6").Value = -34
$ws.Range("T6").Value = 51
$ws.Range("U6").Value = 160
$ws.Range("V6").Value = 1743
$ws.Range("W6").Value = 3.84
$ws.Range("X6").Value = 0.8100000000000001
$ws.Range("Y6").Value = 1.13
$ws.Range("Z6").Value = 0.49
$ws.Range("AA6").Value = 132.54
$ws.Range("AB6").Value = 614.09
$ws.Range("AC6").Value = 46
$ws.Range("AD6").Value = 84.31
$ws.Range("AE6").Value = 4048
$ws.Range("AF6").Value = 0.96
$ws.Range("AG6").Value = 20
$ws.Range("AH6").Value = 0.51
$ws.Range("AI6").Value = 43.4
$ws.Range("AJ6").Value = 38372720

# Row 7
$ws.Range("D7").Value = 2886
$ws.Range("E7").Value = 129
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 70
$ws.Range("I7").Value = 50
$ws.Range("K7").Value = 4054
$ws.Range("L7").Value = 2426
$ws.Range("M7").Value = 1628
$ws.Range("N7").Value = 1594
$ws.Range("P7").Value = 220
$ws.Range("Q7").Value = 12
$ws.Range("R7").Value = -73
$ws.Range("S7").Value = 238
$ws.Range("T7").Value = 50
$ws.Range("U7").Value = -237
$ws.Range("W7").Value = 4.47
$ws.Range("X7").Value = 2.43
$ws.Range("Y7").Value = 3.18
$ws.Range("Z7").Value = 1.81
$ws.Range("AA7").Value = 148.99
$ws.Range("AC7").Value = 130
$ws.Range("AD7").Value = 31.04
$ws.Range("AE7").Value = 4166
$ws.Range("AF7").Value = 0.97
$ws.Range("AG7").Value = 20
$ws.Range("AH7").Value = 0.49
$ws.Range("AI7").Value = 15.35

# Row 8
$ws.Range("D8").Value = 2985
$ws.Range("E8").Value = 169
$ws.Range("G8").Value = 26
$ws.Range("H8").Value = 106
$ws.Range("I8").Value = 82
$ws.Range("K8").Value = 4309
$ws.Range("L8").Value = 2582
$ws.Range("M8").Value = 1728
$ws.Range("N8").Value = 1687
$ws.Range("P8").Value = 220
$ws.Range("Q8").Value = 192
$ws.Range("R8").Value = -101
$ws.Range("S8").Value = 96
$ws.Range("T8").Value = 88
$ws.Range("U8").Value = 120
$ws.Range("W8").Value = 5.66
$ws.Range("X8").Value = 3.57
$ws.Range("Y8").Value = 5.03
$ws.Range("Z8").Value = 2.55
$ws.Range("AA8").Value = 149.44
$ws.Range("AC8").Value = 215
$ws.Range("AD8").Value = 18.81
$ws.Range("AE8").Value = 4436
$ws.Range("AF8").Value = 0.91
$ws.Range("AG8").Value = 20
$ws.Range("AH8").Value = 0.49
$ws.Range("AI8").Value = 9.300000000000001

# Row 9
$ws.Range("D9").Value = 3372
$ws.Range("E9").Value = 210
$ws.Range("G9").Value = 170
$ws.Range("H9").Value = 130
$ws.Range("I9").Value = 122
$ws.Range("K9").Value = 4451
$ws.Range("L9").Value = 2602
$ws.Range("M9").Value = 1850
$ws.Range("N9").Value = 1800
$ws.Range("P9").Value = 220
$ws.Range("Q9").Value = 179
$ws.Range("R9").Value = -110
$ws.Range("S9").Value = -54
$ws.Range("T9").Value = 88
$ws.Range("U9").Value = 109
$ws.Range("W9").Value = 6.23
$ws.Range("X9").Value = 3.84
$ws.Range("Y9").Value = 6.97
$ws.Range("Z9").Value = 2.96
$ws.Range("AA9").Value = 140.66
$ws.Range("AC9").Value = 317
$ws.Range("AD9").Value = 12.78
$ws.Range("AE9").Value = 4734
$ws.Range("AF9").Value = 0.85
$ws.Range("AG9").Value = 20
$ws.Range("AH9").Value = 6.32
